# Update Active_Outages.xlsx - refresh elapsed-duration readouts and
# correct the stray R4 summary row duplicated at the bottom of sheet R2.

$wb = $excel.ActiveWorkbook

# --- Refresh "Elapsed Duration(Hrs)" (column G) on each region sheet ---
$ws = $wb.Worksheets.Item("R1")
$ws.Range("G2").Value = "3959:29:56"
$ws.Range("G3").Value = "99:02:34"
$ws.Range("G4").Value = "122:02:34"

$ws = $wb.Worksheets.Item("R2")
$ws.Range("G2").Value = "12140:53:44"
$ws.Range("G3").Value = "3270:37:13"
$ws.Range("G4").Value = "508:48:47"

# Row 6 on R2 was a stray/duplicate "R4" entry showing a stale Hub Site
# and Battery Backup Status - update it to match the current PCM record.
$ws.Range("D6").Value = "JED0925"
$ws.Range("J6").Value = "In progress"

$ws = $wb.Worksheets.Item("R4")
$ws.Range("G2").Value = "2986:43:26"
$ws.Range("G3").Value = "213:55:41"
$ws.Range("G4").Value = "102:08:06"
$ws.Range("G5").Value = "99:45:39"

$ws = $wb.Worksheets.Item("R5")
$ws.Range("G2").Value = "460:42:25"

$ws = $wb.Worksheets.Item("R6")
$ws.Range("G2").Value = "101:14:43"
